$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# The old sheet had a duplicated "Contact" / "No display for ContactDetail" row
# (rows 10 and 11). Remove the second occurrence (row 11); this shifts every
# row below it up by one.
$ws.Rows.Item(11).Delete()

# --- Simple value updates (rows renumbered after the delete above) ---
$ws.Range("B3").Value = "6.0.0"
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"
$ws.Range("B9").Value = "Alvearie Team"
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# B14 ("Case Sensitive" row) needs the literal text "true" (a shared string),
# not the Boolean TRUE that a plain Value assignment would coerce it to.
# Build it as a text formula result on a scratch cell well outside the used
# range, then paste just the value (not formats) into B14 so the cell keeps
# its original style and the sheet's used range/dimension is untouched.
$scratch = $ws.Range("A100")
$scratch.Formula = "=""true"""
$scratch.Copy()
$ws.Range("B14").PasteSpecial(-4163)
$scratch.Clear()
